$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# get a temporary Text number format so the literal string is preserved,
# then the format is reset back to General/Normal so no visible style diff remains.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.486.46'
$ws.Range('D3').Value = '1.865.97'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  -0.31%  '
Set-TextValue 'D5' '316.27'
$ws.Range('E5').Value = '  +2.57%  '
Set-TextValue 'D6' '1.004'
$ws.Range('E6').Value = '  -0.23%  '
Set-TextValue 'D7' '0.4666'
$ws.Range('E7').Value = '  +1.04%  '
Set-TextValue 'D8' '0.3734'
$ws.Range('E8').Value = '  +2.19%  '
Set-TextValue 'D9' '0.07396'
$ws.Range('E9').Value = '  +2.41%  '
Set-TextValue 'D10' '0.8889'
$ws.Range('E10').Value = '  +3.63%  '
Set-TextValue 'D11' '0.07968'
$ws.Range('E11').Value = '  +5.92%  '
Set-TextValue 'D12' '19.98'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').Value = '1.850.33'
$ws.Range('E13').Value = '  +6.69%  '
$ws.Range('E14').Value = '  +2.11%  '
Set-TextValue 'D15' '6.604'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('E16').Value = '  +1.08%  '
Set-TextValue 'D17' '1.005'
$ws.Range('E17').Value = '  -0.27%  '
Set-TextValue 'D18' '0.000008966'
$ws.Range('E18').Value = '  +4.30%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('E20').Value = '  +3.68%  '
$ws.Range('D21').Value = '27.515.96'
$ws.Range('E21').Value = '  +3.28%  '
Set-TextValue 'D22' '5.179'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = '2.084.04'
$ws.Range('E24').Value = '  +6.50%  '
Set-TextValue 'D25' '153.24'
$ws.Range('E25').Value = '  +0.99%  '
Set-TextValue 'D26' '1.879'
$ws.Range('E26').Value = '  +2.15%  '
Set-TextValue 'D27' '18.54'
$ws.Range('E27').Value = '  +2.41%  '
Set-TextValue 'D28' '2.090'
$ws.Range('E28').Value = '  +0.74%  '
Set-TextValue 'D29' '5.168'
$ws.Range('E29').Value = '  +1.59%  '
Set-TextValue 'D30' '117.18'
$ws.Range('E30').Value = '  +1.84%  '
Set-TextValue 'D31' '0.08897'
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  +1.95%  '
Set-TextValue 'D33' '0.7531'
$ws.Range('E33').Value = '  +5.51%  '
Set-TextValue 'D34' '1.160'
Set-TextValue 'D35' '4.492'
$ws.Range('E35').Value = '  +1.89%  '
Set-TextValue 'D36' '2.644'
$ws.Range('E36').Value = '  +9.75%  '
$ws.Range('E37').Value = '  +2.96%  '
Set-TextValue 'D38' '1.081'
$ws.Range('E38').Value = '  +0.49%  '
Set-TextValue 'D39' '0.05286'
$ws.Range('E39').Value = '  +0.96%  '
Set-TextValue 'D40' '2.986'
$ws.Range('E40').Value = '  +2.28%  '
Set-TextValue 'D41' '7.189'
$ws.Range('E41').Value = '  +0.82%  '
Set-TextValue 'D42' '0.5216'
$ws.Range('E42').Value = '  +1.54%  '
Set-TextValue 'D43' '0.1645'
$ws.Range('E43').Value = '  +1.58%  '
Set-TextValue 'D44' '8.359'
$ws.Range('E44').Value = '  +2.53%  '
$ws.Range('E45').Value = '  +2.19%  '
Set-TextValue 'D46' '10.36'
$ws.Range('E46').Value = '  +2.55%  '
Set-TextValue 'D47' '1.005'
$ws.Range('E47').Value = '  -0.18%  '
Set-TextValue 'D48' '1.668'
$ws.Range('E48').Value = '  +3.35%  '
Set-TextValue 'D49' '103.60'
$ws.Range('E49').Value = '  +0.67%  '
Set-TextValue 'D50' '0.06261'
$ws.Range('E50').Value = '  -0.23%  '
Set-TextValue 'D51' '65.94'
$ws.Range('E51').Value = '  +3.40%  '
